# "Hortaliza, Vega Central Mapocho de Santiago - Rabanito" weekly refresh.
# A new week of price data is inserted; every existing record from row 300
# down shifts one row, and a brand-new row 414 is appended with the record
# that used to be the last row (413). Only the columns that vary per record
# (Fecha, Volumen, Precio minimo/maximo/promedio ponderado, Origen, Precio
# $/Kg) need to move; the rest of the table (market/product metadata) is
# identical on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 413
$newLastRow = 414

# Row 414 does not exist yet: seed it with the constant metadata columns
# shared by every "Rabanito / Vega Central Mapocho de Santiago" record, and
# copy the date column's number format so it matches the rest of column D.
$ws.Range("A$newLastRow").Value = 9
$ws.Range("B$newLastRow").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C$newLastRow").Value = "Metropolitana"
$ws.Range("E$newLastRow").Value = 13
$ws.Range("F$newLastRow").Value = 300000001
$ws.Range("G$newLastRow").Value = "Rabanito"
$ws.Range("H$newLastRow").Value = "Sin especificar"
$ws.Range("I$newLastRow").Value = "Primera"
$ws.Range("N$newLastRow").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("Q$newLastRow").Value = 100
$ws.Range("R$newLastRow").Value = "Hortaliza"
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# New values per row for the columns that change week over week: Fecha (D),
# Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio
# ponderado (M), Origen (O) and Precio $/Kg (P). Row 300 gets the new week's
# figures; rows 301-414 each inherit what used to sit one row above them.
$rowData = @(
    @{Row=300; D=45027; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=301; D=44971; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=302; D=44553; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=303; D=44370; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=304; D=44985; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=305; D=44222; J=15000; K=2500; L=3000; M=2767; O='Provincia de Chacabuco'; P=28}
    @{Row=306; D=44518; J=8800; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=307; D=44420; J=9700; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=308; D=44817; J=7900; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=309; D=44453; J=7900; K=3500; L=4000; M=3747; O='Provincia de Chacabuco'; P=37}
    @{Row=310; D=44750; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=311; D=44757; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=312; D=44463; J=7900; K=3500; L=4000; M=3747; O='Provincia de Chacabuco'; P=37}
    @{Row=313; D=44265; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=314; D=44882; J=13000; K=3000; L=4000; M=3615; O='Provincia de Chacabuco'; P=36}
    @{Row=315; D=44806; J=7900; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=316; D=44953; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=317; D=44301; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=318; D=44616; J=5200; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=319; D=44187; J=19000; K=2500; L=3000; M=2763; O='Provincia de Chacabuco'; P=28}
    @{Row=320; D=44517; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=321; D=44418; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=322; D=44778; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=323; D=44413; J=8800; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=324; D=44266; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=325; D=44469; J=7900; K=3500; L=4000; M=3747; O='Provincia de Chacabuco'; P=37}
    @{Row=326; D=44411; J=6900; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=327; D=44278; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=328; D=44286; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=329; D=44442; J=7900; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=330; D=44950; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=331; D=44791; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=332; D=44362; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=333; D=44397; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=334; D=44988; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=335; D=44635; J=4300; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=336; D=44489; J=6100; K=3000; L=4000; M=3500; O='Provincia de Chacabuco'; P=35}
    @{Row=337; D=45006; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=338; D=44400; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=339; D=44658; J=5200; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=340; D=44699; J=17000; K=2500; L=3000; M=2735; O='Provincia de Chacabuco'; P=27}
    @{Row=341; D=44181; J=14000; K=2500; L=3000; M=2714; O='Provincia de Chacabuco'; P=27}
    @{Row=342; D=44679; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=343; D=44544; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=344; D=44677; J=6100; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=345; D=44516; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=346; D=44259; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=347; D=44615; J=5200; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=348; D=44188; J=12000; K=2500; L=3000; M=2708; O='Provincia de Chacabuco'; P=27}
    @{Row=349; D=44754; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=350; D=44729; J=6100; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=351; D=44537; J=8800; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=352; D=44908; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=353; D=44603; J=6400; K=3000; L=3500; M=3250; O='Provincia de Chacabuco'; P=32}
    @{Row=354; D=44595; J=8000; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=355; D=44189; J=13000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=356; D=44385; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=357; D=44371; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=358; D=45008; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=359; D=44789; J=9700; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=360; D=44873; J=14000; K=3000; L=4000; M=3643; O='Provincia de Chacabuco'; P=36}
    @{Row=361; D=44694; J=6100; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=362; D=44936; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=363; D=44951; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=364; D=44218; J=13000; K=2500; L=3000; M=2692; O='Provincia de Chacabuco'; P=27}
    @{Row=365; D=44512; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=366; D=44910; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=367; D=44343; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=368; D=44232; J=15000; K=2500; L=3000; M=2733; O='Provincia de Chacabuco'; P=27}
    @{Row=369; D=44673; J=5200; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=370; D=44428; J=7900; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=371; D=44638; J=5200; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=372; D=44714; J=5200; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=373; D=44336; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=374; D=44386; J=4000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=375; D=44826; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=376; D=44925; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=377; D=44504; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=378; D=44328; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=379; D=44868; J=11000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=380; D=44687; J=6100; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=381; D=44847; J=6000; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=382; D=44490; J=7900; K=3000; L=4000; M=3494; O='Provincia de Chacabuco'; P=35}
    @{Row=383; D=44987; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=384; D=44777; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=385; D=45021; J=6000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=386; D=44414; J=7900; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=387; D=44824; J=6100; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=388; D=44264; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=389; D=44316; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=390; D=44901; J=9000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=391; D=44176; J=14000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=392; D=44474; J=7900; K=3500; L=4000; M=3747; O='Provincia de Chacabuco'; P=37}
    @{Row=393; D=45002; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=394; D=44705; J=16000; K=2500; L=3000; M=2781; O='Provincia de Chacabuco'; P=28}
    @{Row=395; D=44448; J=7900; K=3500; L=4000; M=3747; O='Provincia de Chacabuco'; P=37}
    @{Row=396; D=45015; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=397; D=44342; J=9000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=398; D=44186; J=4000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=399; D=44441; J=7900; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=400; D=44952; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=401; D=44727; J=6000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=402; D=44252; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=403; D=44855; J=5000; K=4000; L=4000; M=4000; O='Provincia de Chacabuco'; P=40}
    @{Row=404; D=44391; J=4300; K=2500; L=3000; M=2750; O='Región Metropolitana'; P=28}
    @{Row=405; D=44168; J=21000; K=2500; L=3000; M=2738; O='Provincia de Chacabuco'; P=27}
    @{Row=406; D=44642; J=4300; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=407; D=44454; J=6100; K=3500; L=4000; M=3750; O='Provincia de Chacabuco'; P=38}
    @{Row=408; D=44426; J=6100; K=3000; L=3500; M=3250; O='Provincia de Chacabuco'; P=32}
    @{Row=409; D=44526; J=7900; K=2500; L=3000; M=2747; O='Provincia de Chacabuco'; P=27}
    @{Row=410; D=44551; J=6100; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=411; D=44918; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=412; D=44243; J=10000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=413; D=44217; J=11000; K=2500; L=3000; M=2727; O='Provincia de Chacabuco'; P=27}
    @{Row=414; D=45007; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
)

foreach ($rec in $rowData) {
    $r = $rec.Row
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("J$r").Value = $rec.J
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
}
